$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (Total) sheet,
#    matching the layout/formatting of the other quarterly fund sheets.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Copy formatting (bold header row + bordered style, bold column A) from
# the template sheet so the new sheet matches the existing visual style.
$template.Range("A1:H8").Copy()
$newSheet.Range("A1:H8").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A9").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings data for 2022-Q1 - force text format so codes keep
# leading zeros and decimal-looking figures aren't coerced to numbers.
$newSheet.Range("B2:G9").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "513090"
$newSheet.Range("C2").Value = "易方达中证香港证券投资主题ETF"
$newSheet.Range("D2").Value = "11.07"
$newSheet.Range("E2").Value = "96.47"
$newSheet.Range("F2").Value = "4.03"
$newSheet.Range("G2").Value = "0.4461"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011722"
$newSheet.Range("C3").Value = "前海开源深圳特区精选股票型证券投资基金A"
$newSheet.Range("D3").Value = "3.00"
$newSheet.Range("E3").Value = "85.40"
$newSheet.Range("F3").Value = "6.76"
$newSheet.Range("G3").Value = "0.2028"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "004496"
$newSheet.Range("C4").Value = "前海开源多元策略灵活配置混合A"
$newSheet.Range("D4").Value = "3.09"
$newSheet.Range("E4").Value = "91.25"
$newSheet.Range("F4").Value = "5.97"
$newSheet.Range("G4").Value = "0.1845"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "004497"
$newSheet.Range("C5").Value = "前海开源多元策略灵活配置混合C"
$newSheet.Range("D5").Value = "2.06"
$newSheet.Range("E5").Value = "91.25"
$newSheet.Range("F5").Value = "5.97"
$newSheet.Range("G5").Value = "0.1230"
$newSheet.Range("H5").Value = 8

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "011723"
$newSheet.Range("C6").Value = "前海开源深圳特区精选股票型证券投资基金C"
$newSheet.Range("D6").Value = "0.44"
$newSheet.Range("E6").Value = "85.40"
$newSheet.Range("F6").Value = "6.76"
$newSheet.Range("G6").Value = "0.0297"
$newSheet.Range("H6").Value = 7

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "161124"
$newSheet.Range("C7").Value = "易方达香港恒生综合小型股指数（QDII-LOF）A"
$newSheet.Range("D7").Value = "0.28"
$newSheet.Range("E7").Value = "92.62"
$newSheet.Range("F7").Value = "1.33"
$newSheet.Range("G7").Value = "0.0037"
$newSheet.Range("H7").Value = 8

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "006263"
$newSheet.Range("C8").Value = "易方达香港恒生综合小型股指数（QDII-LOF）C"
$newSheet.Range("D8").Value = "0.06"
$newSheet.Range("E8").Value = "92.62"
$newSheet.Range("F8").Value = "1.33"
$newSheet.Range("G8").Value = "0.0008"
$newSheet.Range("H8").Value = 8

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "002860"
$newSheet.Range("C9").Value = "前海开源沪港深新机遇灵活配置混合"
$newSheet.Range("D9").Value = "0.01"
$newSheet.Range("E9").Value = "83.26"
$newSheet.Range("F9").Value = "6.14"
$newSheet.Range("G9").Value = "0.0006"
$newSheet.Range("H9").Value = 9

# ---------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row at the top of the "总计" sheet's data,
#    pushing the existing quarters down by one row.
#    NOTE: re-resolve the "总计" worksheet reference now that a new sheet
#    has been inserted in front of it - sheet handles in this engine are
#    positional, so the original $totalSheet variable would otherwise now
#    point at the freshly-added "2022-Q1" sheet instead.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 0.99

# The row-insert shifted the old index column (A) down without updating
# its values, so the running 0-based counter needs to be bumped by one
# for every pre-existing quarter row.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
